# Apply edit: add two new worksheets ("PrototypeOneTest output" and
# "PrototypeTwoTest output") containing sample console output, matching the
# commit "Added sample output to CategoryPrices excel under src/test/resources".

$wb = $excel.ActiveWorkbook

# --- Locate existing sheets ---
$tables = $wb.Worksheets.Item("Tables")
$manualCheck = $wb.Worksheets.Item("ManualCheck")

# --- Add the two new worksheets after "ManualCheck" ---
# "PrototypeTwoTest output" is created first (so it receives the lower
# sheetId=4), then "PrototypeOneTest output" is inserted immediately before
# it (so it receives sheetId=5), matching the target tab order:
#   Tables, ManualCheck, PrototypeOneTest output, PrototypeTwoTest output
$newTwo = $wb.Worksheets.Add($null, $manualCheck)
$newTwo.Name = "PrototypeTwoTest output"

$newOne = $wb.Worksheets.Add($wb.Worksheets.Item("PrototypeTwoTest output"))
$newOne.Name = "PrototypeOneTest output"

# Re-fetch fresh, stable references to the new sheets now that both exist
# and have their final names/positions (avoids stale COM references).
$wsOne = $wb.Worksheets.Item("PrototypeOneTest output")
$wsTwo = $wb.Worksheets.Item("PrototypeTwoTest output")

# --- Content for "PrototypeOneTest output" (column A, rows 1-26) ---
$sheet3Content = @(
    'Enter Newspaper whose monthly subscription needs to be calculated. Choices are :: [Hindu, Carwash, BM, HT, Magazine, TOI, ET, Milk]',
    'The input to the program should be a comma separated list of papers that the consumer wants to subscribe (eg: TOI,ET)',
    'Enter your string here :: TOI,BM,HT,Random,TOI,Random,BM,@#@$,BM',
    '[BM, HT, TOI]',
    $null,
    '*******************************************************************',
    $null,
    'Enter the Year followed by Month (using first three letters only) seperated by comma whose budget needs calculated (eg: 2020,Mar) :: 2020,Apr',
    'Total number of days in the Month ''Apr'' and Year ''2020'' is :: 30',
    '{Monday=4, Thursday=5, Friday=4, Sunday=4, Wednesday=5, Tuesday=4, Saturday=4}',
    $null,
    '*******************************************************************',
    $null,
    'Calculation Table for ''BM'' is :: {Monday=1.5, Thursday=1.5, Friday=1.5, Sunday=1.5, Wednesday=1.5, Tuesday=1.5, Saturday=1.5}',
    'Total Monthly subscription for BM is Rs 45.0',
    $null,
    'Calculation Table for ''HT'' is :: {Monday=2.0, Thursday=2.0, Friday=2.0, Sunday=4.0, Wednesday=2.0, Tuesday=2.0, Saturday=4.0}',
    'Total Monthly subscription for HT is Rs 76.0',
    $null,
    'Calculation Table for ''TOI'' is :: {Monday=3.0, Thursday=3.0, Friday=3.0, Sunday=6.0, Wednesday=3.0, Tuesday=3.0, Saturday=5.0}',
    'Total Monthly subscription for TOI is Rs 110.0',
    $null,
    $null,
    '*******************************************************************',
    $null,
    'Total estimated subscription amount for the month of ''Apr'' and year ''2020'' is  :: Rs 231.0',
)

for ($i = 0; $i -lt $sheet3Content.Length; $i++) {
    $val = $sheet3Content[$i]
    if ($null -ne $val) {
        $wsOne.Cells.Item($i + 1, 1).Value = $val
    }
}

# --- Content for "PrototypeTwoTest output" (column A, rows 1-36) ---
$sheet4Content = @(
    'Enter category followed by frequency whose subscription needs to be calculated. Choices are :: [Hindu, Carwash, BM, HT, Magazine, TOI, ET, Milk]',
    'The input to the program should be a comma separated list of categories that the consumer wants to subscribe followed by Daily,Weekly,Biweekly (eg: TOI,Daily,ET,Weekly,HM,BiWeekly)',
    'Enter your string here :: BM,Daily,BM,Daily,ET,Blah,Blah,Weekly,HT,Daily,TOI,Daily,Milk,Daily,Magazine,Weekly,Carwash,Biweekly',
    '{Carwash=Biweekly, BM=Daily, HT=Daily, Magazine=Weekly, TOI=Daily, Milk=Daily}',
    '*******************************************************************',
    $null,
    'Enter the Year followed by Month (using first three letters only) seperated by comma whose budget needs calculated (eg: 2020,Mar) :: 2020,Apr',
    'Total number of days in the Month ''Apr'' and Year ''2020'' is :: 30',
    '{Monday=4, Thursday=5, Friday=4, Sunday=4, Wednesday=5, Tuesday=4, Saturday=4}',
    $null,
    '*******************************************************************',
    $null,
    'Calculation Table for ''Carwash'' is :: {Monday=200.0, Thursday=200.0, Friday=200.0, Sunday=200.0, Wednesday=200.0, Tuesday=200.0, Saturday=200.0}',
    'Total Biweekly subscription for Carwash is Rs 200.0',
    'Therefore, total monthly subscription for Carwash is Rs 400.0',
    $null,
    'Calculation Table for ''BM'' is :: {Monday=1.5, Thursday=1.5, Friday=1.5, Sunday=1.5, Wednesday=1.5, Tuesday=1.5, Saturday=1.5}',
    'Total Monthly subscription for BM is Rs 45.0',
    $null,
    'Calculation Table for ''HT'' is :: {Monday=2.0, Thursday=2.0, Friday=2.0, Sunday=4.0, Wednesday=2.0, Tuesday=2.0, Saturday=4.0}',
    'Total Monthly subscription for HT is Rs 76.0',
    $null,
    'Calculation Table for ''Magazine'' is :: {Monday=150.0, Thursday=150.0, Friday=150.0, Sunday=150.0, Wednesday=150.0, Tuesday=150.0, Saturday=150.0}',
    'Total Weekly subscription for Magazine is Rs 150.0',
    'Therefore, total Monthly subscription for Magazine is Rs 600.0',
    $null,
    'Calculation Table for ''TOI'' is :: {Monday=3.0, Thursday=3.0, Friday=3.0, Sunday=6.0, Wednesday=3.0, Tuesday=3.0, Saturday=5.0}',
    'Total Monthly subscription for TOI is Rs 110.0',
    $null,
    'Calculation Table for ''Milk'' is :: {Monday=50.0, Thursday=50.0, Friday=50.0, Sunday=50.0, Wednesday=50.0, Tuesday=50.0, Saturday=50.0}',
    'Total Monthly subscription for Milk is Rs 1500.0',
    $null,
    $null,
    '*******************************************************************',
    $null,
    'Total estimated subscription amount for the month of ''Apr'' and year ''2020'' is  :: Rs 2731.0',
)

for ($i = 0; $i -lt $sheet4Content.Length; $i++) {
    $val = $sheet4Content[$i]
    if ($null -ne $val) {
        $wsTwo.Cells.Item($i + 1, 1).Value = $val
    }
}

# --- Keep "Tables" as the active/selected sheet ---
$tables.Select()

# --- Cosmetic workbook view setting from the diff (best-effort; engine may
# not persist this attribute in the OOXML, but we set it for fidelity) ---
$wb.Windows.Item(1).TabRatio = 0.813

Write-Host "Sheets after edit:"
foreach ($s in $wb.Worksheets) {
    Write-Host (" - " + $s.Name)
}
